$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump the generation timestamp ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-10-29T11:46:56+00:00"

# --- Elements sheet: fix casing / wording of the ExerciceProfessionnel row ---
$elements = $wb.Worksheets.Item("Elements")

# ID + Path columns (A & B): exerciceProfessionnel -> ExerciceProfessionnel
$elements.Range("A6").Value = "Competence.ExerciceProfessionnel"
$elements.Range("B6").Value = "Competence.ExerciceProfessionnel"

# Short + Definition columns (L & M): drop the trailing period
$elements.Range("L6").Value = "Lien vers la classe ExerciceProfessionnel"
$elements.Range("M6").Value = "Lien vers la classe ExerciceProfessionnel"

# Base Path column (AF): exerciceProfessionnel -> ExerciceProfessionnel
$elements.Range("AF6").Value = "SavoirFaire.ExerciceProfessionnel"
